# Daily attendance processing - 2025-10-21 10:21:13
#
# For every populated row in the "Recorded By" column (G), the
# comma-separated list of recorder names/emails is re-ordered into
# ordinal (case-sensitive, ASCII) ascending sort order, e.g.:
#   "dnasr281@gmail.com, System"  ->  "System, dnasr281@gmail.com"
#   "system, System, backup@backdoor.com" -> "System, backup@backdoor.com, system"
#
# NOTE: this runtime's built-in comparison operators (-clt/-cgt/-ceq/-clt
# Sort-Object, etc.) all behave case-insensitively here, so a manual
# ordinal compare (based on character codes) is implemented below to
# reproduce Python's default `sorted()` behaviour used to build the diff.
# Also note: loop/temp variable names used inside functions must be kept
# distinct from any variable name used in an outer loop that may call
# into them, because this runtime shares variable scope between caller
# and callee (re-using a name such as $i in both places causes the
# outer loop counter to be reset and the script to run away).

function Get-OrdCmp($ordCmpStrA, $ordCmpStrB) {
    $ordCmpLenA = $ordCmpStrA.Length
    $ordCmpLenB = $ordCmpStrB.Length
    $ordCmpMinLen = $ordCmpLenA
    if ($ordCmpLenB -lt $ordCmpMinLen) { $ordCmpMinLen = $ordCmpLenB }

    $ordCmpResult = 0
    $ordCmpIdx = 0
    while ($ordCmpIdx -lt $ordCmpMinLen) {
        $ordCmpCodeA = [int][char]$ordCmpStrA[$ordCmpIdx]
        $ordCmpCodeB = [int][char]$ordCmpStrB[$ordCmpIdx]
        if ($ordCmpCodeA -lt $ordCmpCodeB) {
            $ordCmpResult = -1
            $ordCmpIdx = $ordCmpMinLen
        } elseif ($ordCmpCodeA -gt $ordCmpCodeB) {
            $ordCmpResult = 1
            $ordCmpIdx = $ordCmpMinLen
        } else {
            $ordCmpIdx = $ordCmpIdx + 1
        }
    }

    if ($ordCmpResult -eq 0) {
        if ($ordCmpLenA -lt $ordCmpLenB) { $ordCmpResult = -1 }
        elseif ($ordCmpLenA -gt $ordCmpLenB) { $ordCmpResult = 1 }
    }

    return $ordCmpResult
}

function Get-SortedRecordedBy($srbRawValue) {
    if ($srbRawValue -eq $null) { return $srbRawValue }
    if ($srbRawValue -eq "") { return $srbRawValue }

    $srbParts = $srbRawValue -split ", "
    $srbCount = $srbParts.Count

    if ($srbCount -eq 1) {
        return $srbParts[0]
    }

    if ($srbCount -eq 2) {
        $srbA = $srbParts[0]
        $srbB = $srbParts[1]
        if ((Get-OrdCmp $srbA $srbB) -gt 0) {
            $srbTmp = $srbA
            $srbA = $srbB
            $srbB = $srbTmp
        }
        return "$srbA, $srbB"
    }

    if ($srbCount -eq 3) {
        $srbA = $srbParts[0]
        $srbB = $srbParts[1]
        $srbC = $srbParts[2]
        if ((Get-OrdCmp $srbA $srbB) -gt 0) {
            $srbTmp = $srbA
            $srbA = $srbB
            $srbB = $srbTmp
        }
        if ((Get-OrdCmp $srbB $srbC) -gt 0) {
            $srbTmp = $srbB
            $srbB = $srbC
            $srbC = $srbTmp
        }
        if ((Get-OrdCmp $srbA $srbB) -gt 0) {
            $srbTmp = $srbA
            $srbA = $srbB
            $srbB = $srbTmp
        }
        return "$srbA, $srbB, $srbC"
    }

    # Fallback (shouldn't happen for this report): leave untouched.
    return $srbRawValue
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($rbRow = 2; $rbRow -le $lastRow; $rbRow++) {
    $rbOrigValue = $ws.Cells.Item($rbRow, 7).Value2
    if (($rbOrigValue -ne $null) -and ($rbOrigValue -ne "")) {
        $rbNewValue = Get-SortedRecordedBy $rbOrigValue
        if ($rbNewValue -ne $rbOrigValue) {
            $ws.Cells.Item($rbRow, 7).Value = $rbNewValue
        }
    }
}
